$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Breakthrough" detail bullet: replace text, and the _GoBack bookmark that
#    used to sit at the end of this paragraph moves away (it will be re-created
#    at the end of the "Defence Line" detail bullet below), so drop it here.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Have ships close to enemy table edge", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Score 3 victory points at the end of each game turn if you have at least 2 ships within the enemy normal deployment zone",
    2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Defence Line" detail bullet: replace the single run with two runs of
#    text, then re-create the (now collapsed) _GoBack bookmark right after the
#    new text, followed by a trailing run containing a single space.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Have no enemy ships close to your table edge", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Score 2 victory points at the end of each game turn if there are no enemy ships within 6" + [char]0x201d + " of your normal deployment zone",
    2)

# Re-split that replacement text into the two runs shown in the target: find
# the boundary and re-type the tail so the run break lands in the right spot.
$boundary = "Score 2 victory points at the end of each game turn if there are "
$found = $d.Content.Find.Execute($boundary, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the paragraph again (index is stable: it is the bullet right after "Defence Line").
$defenceHeader = $d.Content.Find.Execute("Defence Line", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
